# Add a "Login Details" sheet (SauceDemo login credentials) after the
# existing "User Details" sheet, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $usersSheet)
$ws.Name = "Login Details"

# Header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Usernames column, then passwords column (matches original authoring order)
$ws.Range("A2").Value = "standard_user"
$ws.Range("A3").Value = "locked_out_user"
$ws.Range("A4").Value = "problem_user"
$ws.Range("A5").Value = "performance_glitch_user"

$ws.Range("B2").Value = "secret_sauce"
$ws.Range("B3").Value = "secret_sauce"
$ws.Range("B4").Value = "secret_sauce"
$ws.Range("B5").Value = "secret_sauce"

# A couple of the password cells picked up an explicit (non-bold) font
# when they were filled in
$ws.Range("B2:B3").Font.Name = "Calibri"

# Column widths to fit the new content
$ws.Columns.Item(1).ColumnWidth = 23.21875
$ws.Columns.Item(2).ColumnWidth = 16.109375

# Selection / active cell on the new sheet
[void]$ws.Range("B10").Select()

[void]$ws.Activate()
